# Add descriptive metadata (Value column) for several L2130i variables and
# restyle the last three added cells with an explicit black font color.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new "Value" column entries in the order the new shared strings were
# introduced (var.baseline_shift, var.slope_shift, var.residuals, var.H2O,
# var.CH4, var.OutletValve, var.ValveMask, var.WarmBoxTemp, var.DasTemp,
# var.EtalonTemp, var.INST_STATUS).
$ws.Range("B38").Value = "change in constant term of fitted baseline relative to empty cavity baseline; see Rev C 04-2015 manual"
$ws.Range("B39").Value = "change in linear term of fitted baseline relative to empty cavity baseline; see Rev C 04-2015 manual"
$ws.Range("B40").Value = "MRS residual of LS fit of measured spectra versus the expected spectra; see Rev C 04-2015 manual"

$ws.Range("B33").Value = "Humidity; ppmv"
$ws.Range("B36").Value = "Methane concentration in ppm"

$ws.Range("B31").Value = 'Relative value of degree of "openness" of cavity valve; continuously variable'
$ws.Range("B32").Value = "Indicates which valves in vaporizer are activated; convert to binary to determine which valves are open"

$ws.Range("B27").Value = "Temperature at warm box in instrument interior"
$ws.Range("B28").Value = "Temperature sensor in the interior of the analyzer"
$ws.Range("B29").Value = "Temperature sensor in the interior of the analyzer"

$ws.Range("B24").Value = "Presumably code that states whether the instrument is operating as expected; not documented anywhere by Picarro"

# Give the last three newly-populated cells (the "Rev C" baseline/slope/
# residual descriptions) an explicit black font color, introducing a new
# font + cell style.
$ws.Range("B38").Font.Color = 0
$ws.Range("B39").Font.Color = 0
$ws.Range("B40").Font.Color = 0

# Move the active selection to B15 (matches the updated sheet view state).
$ws.Range("B15").Select() | Out-Null
